$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between F and L columns for rows 4 and 8
$ws.Range("F4").Value = 0
$ws.Range("L4").Value = "-"

$ws.Range("F8").Value = 0
$ws.Range("L8").Value = "-"
